# Remove the "By <Name>" attribution paragraph that follows the slide
# title on slides 3-14 (1-based Slides.Item index). Each of those title
# placeholders holds exactly two paragraphs: the real title, then a
# second paragraph whose only run is "By <Name>".
#
# Simply deleting the trailing paragraph (TextRange.Paragraphs(2,1).Delete())
# removes its text but the serializer still leaves behind an empty
# trailing <a:p/> stub, because it is the *last* paragraph of the text
# frame. To get a byte-for-byte removal of the whole <a:p> element we
# instead:
#   1. copy the first paragraph's text (title) into the second
#      paragraph (so the trailing/"last" paragraph becomes the one we
#      want to keep, formatting included),
#   2. delete the now-redundant first paragraph, which is no longer the
#      last paragraph, so it is removed cleanly with no stub left over.

$p = $ppt.ActivePresentation

for ($i = 3; $i -le 14; $i++) {
    $slide = $p.Slides.Item($i)
    $titleShape = $slide.Shapes.Item(1)
    $tr = $titleShape.TextFrame.TextRange

    if ($tr.Paragraphs().Count -lt 2) {
        continue
    }

    $firstPara = $tr.Paragraphs(1, 1)
    # Paragraph .Text includes the trailing paragraph-mark (CR); strip it
    # before reusing the string as plain run text.
    $titleText = $firstPara.Text.TrimEnd([char]13)

    $secondPara = $tr.Paragraphs(2, 1)
    # Clear the paragraph's text first so the engine's diff-based text
    # setter can't preserve a run fragment shared between the old
    # ("By <Name>") and new (title) text as a leftover common-suffix
    # run (e.g. "...Risk" / "...Malak" both ending in "k") - clearing
    # first guarantees a single clean run with no stray <a:rPr>.
    $secondPara.Text = ""
    $secondPara = $tr.Paragraphs(2, 1)
    $secondPara.Text = $titleText

    # Re-fetch paragraph 1 (still the "By <Name>" formatting-less title
    # duplicate) and delete it; since it is no longer the last paragraph
    # in the text frame it is removed cleanly, with no leftover <a:p/>.
    $firstParaAgain = $tr.Paragraphs(1, 1)
    $firstParaAgain.Delete()
}
